$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The new template has 16 columns (A:P) instead of 17 (A:Q); drop the old trailing column.
$ws.Columns.Item(17).Delete()

# Columns A and B now hold date-looking text ("발주일자"/"납기일자"); force Text format first so
# Excel does not silently convert the strings to date serial numbers.
$ws.Range("A2:B6").NumberFormat = "@"

# Row 1
$ws.Range("A1").Value = "발주일자"
$ws.Range("B1").Value = "납기일자"
$ws.Range("C1").Value = "거래처명"
$ws.Range("D1").Value = "거래처 이메일"
$ws.Range("E1").Value = "납품처명"
$ws.Range("F1").Value = "납품처 이메일"
$ws.Range("G1").Value = "프로젝트명"
$ws.Range("H1").Value = "대분류"
$ws.Range("I1").Value = "중분류"
$ws.Range("J1").Value = "소분류"
$ws.Range("K1").Value = "품목명"
$ws.Range("L1").Value = "규격"
$ws.Range("M1").Value = "수량"
$ws.Range("N1").Value = "단가"
$ws.Range("O1").Value = "총금액"
$ws.Range("P1").Value = "비고"

# Row 2
$ws.Range("A2").Value = "2025-08-25"
$ws.Range("B2").Value = "2025-09-06"
$ws.Range("C2").Value = "신호수"
$ws.Range("D2").Value = "신호수@example.com"
$ws.Range("E2").Value = "힐스테이트 도곡동1차"
$ws.Range("F2").Value = "delivery@example.com"
$ws.Range("G2").Value = "힐스테이트 도곡동1차"
$ws.Range("H2").Value = "4. 장비비"
$ws.Range("I2").Value = "2) 신호수 외"
$ws.Range("J2").Value = "기타"
$ws.Range("K2").Value = "4월"
$ws.Range("L2").Value = "KS규격-1"
$ws.Range("M2").Value = 2
$ws.Range("N2").Value = 155000
$ws.Range("O2").Value = 341000
$ws.Range("P2").Value = "탁영롱"

# Row 3
$ws.Range("A3").Value = "2025-08-31"
$ws.Range("B3").Value = "2025-09-11"
$ws.Range("C3").Value = "신호수"
$ws.Range("D3").Value = "신호수@example.com"
$ws.Range("E3").Value = "힐스테이트 도곡동1차"
$ws.Range("F3").Value = "delivery@example.com"
$ws.Range("G3").Value = "힐스테이트 도곡동1차"
$ws.Range("H3").Value = "4. 장비비"
$ws.Range("I3").Value = "2) 신호수 외"
$ws.Range("J3").Value = "기타"
$ws.Range("K3").Value = "5월"
$ws.Range("L3").Value = "KS규격-2"
$ws.Range("M3").Value = 6
$ws.Range("N3").Value = 155000
$ws.Range("O3").Value = 1023000
$ws.Range("P3").Value = "탁영롱 "

# Row 4
$ws.Range("A4").Value = "2025-08-27"
$ws.Range("B4").Value = "2025-10-16"
$ws.Range("C4").Value = "신호수"
$ws.Range("D4").Value = "신호수@example.com"
$ws.Range("E4").Value = "힐스테이트 도곡동1차"
$ws.Range("F4").Value = "delivery@example.com"
$ws.Range("G4").Value = "힐스테이트 도곡동1차"
$ws.Range("H4").Value = "4. 장비비"
$ws.Range("I4").Value = "2) 신호수 외"
$ws.Range("J4").Value = "기타"
$ws.Range("K4").Value = "4월"
$ws.Range("L4").Value = "KS규격-3"
$ws.Range("M4").Value = 2
$ws.Range("N4").Value = 155000
$ws.Range("O4").Value = 341000
$ws.Range("P4").Value = "서진원"

# Row 5
$ws.Range("A5").Value = "2025-08-23"
$ws.Range("B5").Value = "2025-09-28"
$ws.Range("C5").Value = "신호수"
$ws.Range("D5").Value = "신호수@example.com"
$ws.Range("E5").Value = "힐스테이트 도곡동1차"
$ws.Range("F5").Value = "delivery@example.com"
$ws.Range("G5").Value = "힐스테이트 도곡동1차"
$ws.Range("H5").Value = "4. 장비비"
$ws.Range("I5").Value = "2) 신호수 외"
$ws.Range("J5").Value = "기타"
$ws.Range("K5").Value = "5월"
$ws.Range("L5").Value = "KS규격-4"
$ws.Range("M5").Value = 2.1
$ws.Range("N5").Value = 155000
$ws.Range("O5").Value = 358050
$ws.Range("P5").Value = "손명산 "

# Row 6
$ws.Range("A6").Value = "2025-09-16"
$ws.Range("B6").Value = "2025-09-22"
$ws.Range("C6").Value = "신호수"
$ws.Range("D6").Value = "신호수@example.com"
$ws.Range("E6").Value = "힐스테이트 도곡동1차"
$ws.Range("F6").Value = "delivery@example.com"
$ws.Range("G6").Value = "힐스테이트 도곡동1차"
$ws.Range("H6").Value = "4. 장비비"
$ws.Range("I6").Value = "2) 신호수 외"
$ws.Range("J6").Value = "기타"
$ws.Range("K6").Value = "5월"
$ws.Range("L6").Value = "KS규격-5"
$ws.Range("M6").Value = 1
$ws.Range("N6").Value = 155000
$ws.Range("O6").Value = 170500
$ws.Range("P6").Value = "김병호 "

# Restore the plain (unstyled) look: the template header row is no longer bold/bordered,
# and the date columns should not keep the temporary Text number format either.
$ws.Range("A1:P1").Style = "Normal"
$ws.Range("A2:B6").Style = "Normal"
